# Re-rank the classification table: rows are re-sorted by (new) lap time,
# drivers/teams/numbers reshuffled accordingly, and lap times updated to the
# post-pitstop-ban qualifying times.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Fernando Alonso"
$ws.Range("B2").Value = "Spain"
$ws.Range("C2").Value = "Renault"
$ws.Range("D2").Value = 14
$ws.Range("F2").Value = "0:01:28.926388"

$ws.Range("A3").Value = "Lewis Hamilton"
$ws.Range("B3").Value = "England"
$ws.Range("C3").Value = "Mercedes"
$ws.Range("D3").Value = 44
$ws.Range("F3").Value = "0:01:28.997629"

$ws.Range("A4").Value = "Carlos Sainz Jr"
$ws.Range("B4").Value = "Spain"
$ws.Range("C4").Value = "Ferrari"
$ws.Range("D4").Value = 55
$ws.Range("F4").Value = "0:01:29.044263"

$ws.Range("A5").Value = "Max Verstappen"
$ws.Range("B5").Value = "Netherlands"
$ws.Range("C5").Value = "Red Bull"
$ws.Range("D5").Value = 33
$ws.Range("F5").Value = "0:01:29.047641"

$ws.Range("A6").Value = "Mick Schumacher"
$ws.Range("B6").Value = "Germany"
$ws.Range("C6").Value = "Haas"
$ws.Range("D6").Value = 47
$ws.Range("F6").Value = "0:01:29.191375"

$ws.Range("A7").Value = "Esteban Ocon"
$ws.Range("B7").Value = "France"
$ws.Range("C7").Value = "Renault"
$ws.Range("D7").Value = 31
$ws.Range("F7").Value = "0:01:29.222135"

$ws.Range("A8").Value = "Daniel Ricciardo"
$ws.Range("B8").Value = "Australia"
$ws.Range("C8").Value = "McLaren"
$ws.Range("D8").Value = 3
$ws.Range("F8").Value = "0:01:29.361194"

$ws.Range("A9").Value = "Charles Leclerc"
$ws.Range("B9").Value = "Monaco"
$ws.Range("C9").Value = "Ferrari"
$ws.Range("D9").Value = 16
$ws.Range("F9").Value = "0:01:29.368085"

$ws.Range("A10").Value = "Lando Norris"
$ws.Range("B10").Value = "England"
$ws.Range("C10").Value = "McLaren"
$ws.Range("D10").Value = 4
$ws.Range("F10").Value = "0:01:29.373747"

$ws.Range("A11").Value = "Valteri Bottas"
$ws.Range("B11").Value = "Finland"
$ws.Range("C11").Value = "Mercedes"
$ws.Range("D11").Value = 77
$ws.Range("F11").Value = "0:01:29.439958"

$ws.Range("A12").Value = "Lance Stroll"
$ws.Range("B12").Value = "Canada"
$ws.Range("C12").Value = "Aston Martin"
$ws.Range("D12").Value = 18
$ws.Range("F12").Value = "0:01:29.462654"

$ws.Range("A13").Value = "Nikita Mazepin"
$ws.Range("B13").Value = "Neutral"
$ws.Range("C13").Value = "Haas"
$ws.Range("D13").Value = 9
$ws.Range("F13").Value = "0:01:29.464090"

$ws.Range("A14").Value = "Yuki Tsunoda"
$ws.Range("B14").Value = "Japan"
$ws.Range("C14").Value = "Alpha Tauri"
$ws.Range("D14").Value = 22
$ws.Range("F14").Value = "0:01:29.507213"

$ws.Range("A15").Value = "Antonio Giovinazzi"
$ws.Range("B15").Value = "Italy"
$ws.Range("C15").Value = "Alfa-Romeo"
$ws.Range("D15").Value = 99
$ws.Range("F15").Value = "0:01:29.523189"

$ws.Range("A16").Value = "Pierre Gasly"
$ws.Range("B16").Value = "France"
$ws.Range("C16").Value = "Alpha Tauri"
$ws.Range("D16").Value = 10
$ws.Range("F16").Value = "0:01:29.526840"

$ws.Range("A17").Value = "Sergio Pérez"
$ws.Range("B17").Value = "Mexico"
$ws.Range("C17").Value = "Red Bull"
$ws.Range("D17").Value = 11
$ws.Range("F17").Value = "0:01:29.544101"

$ws.Range("A18").Value = "Kimi Raikkonen"
$ws.Range("B18").Value = "Finland"
$ws.Range("C18").Value = "Alfa-Romeo"
$ws.Range("D18").Value = 7
$ws.Range("F18").Value = "0:01:29.577895"

$ws.Range("A19").Value = "Nicholas Latifi"
$ws.Range("B19").Value = "Canada"
$ws.Range("C19").Value = "Williams"
$ws.Range("D19").Value = 6
$ws.Range("F19").Value = "0:01:30.149072"

$ws.Range("A20").Value = "George Russel"
$ws.Range("B20").Value = "England"
$ws.Range("C20").Value = "Williams"
$ws.Range("D20").Value = 63
$ws.Range("F20").Value = "0:01:30.188002"

$ws.Range("A21").Value = "Sebastian Vettel"
$ws.Range("B21").Value = "Germany"
$ws.Range("C21").Value = "Aston Martin"
$ws.Range("D21").Value = 5
$ws.Range("F21").Value = "0:01:30.361087"

